$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for B1:L19 (columns B..L), one array per row (rows 1..19)
$data = @(
    @(45627, 45634, 45641, 45648, 45655, 45662, 45669, 45676, 45683, 45690, 45697),
    @(64, 61, 58.9, 57.8, 57.8, 58.7, 60.5, 62.9, 65.5, 68, 70.1),
    @(70.1, 65.7, 63, 62, 62.6, 64.7, 67.9, 71.8, 75.7, 79.2, 81.7),
    @(47.2, 45.6, 44.9, 45.2, 46.7, 49.2, 52.8, 57.2, 62, 67, 71.5),
    @(58.5, 56.1, 54.7, 54.3, 54.9, 56.3, 58.3, 60.6, 63, 65.2, 66.8),
    @(68.7, 65.9, 64, 63.3, 63.7, 65.1, 67.1, 69.3, 71.4, 72.8, 73.4),
    @(64.3, 57.7, 51.4, 45.8, 41.1, 37.6, 35.2, 34, 33.8, 34.5, 35.7),
    @(66.2, 62.7, 60.4, 59.4, 59.7, 61, 63.3, 66, 68.9, 71.3, 73.1),
    @(72.3, 66.8, 62.5, 59.6, 58.2, 58.3, 59.8, 62.3, 65.5, 68.9, 72.1),
    @(53.9, 51.6, 49.9, 49.3, 49.8, 51.3, 53.9, 57.4, 61.4, 65.5, 69.3),
    @(60, 60.1, 60.7, 61.9, 63.4, 65.1, 66.8, 68.3, 69.2, 69.4, 68.8),
    @(61.5, 59, 57.6, 57.3, 58.4, 60.6, 63.8, 67.7, 71.9, 75.9, 79.3),
    @(74.5, 70.7, 68.8, 68.6, 70.1, 73, 76.8, 80.8, 84.6, 87.5, 89.1),
    @(60.5, 60.2, 59.6, 58.8, 57.8, 56.5, 54.9, 52.9, 50.6, 48, 45.1),
    @(62.1, 59.4, 57.8, 57.4, 58.3, 60.2, 63, 66.3, 69.7, 72.8, 75.2),
    @(67.4, 63.5, 60.3, 58.1, 57.1, 57.3, 58.5, 60.7, 63.5, 66.4, 69.1),
    @(64, 61.9, 59.8, 57.7, 55.8, 54, 52.2, 50.6, 49, 47.2, 45.4),
    @(69.5, 67.1, 65.5, 64.9, 65.8, 68, 71.7, 76.4, 81.9, 87.5, 92.8),
    @(67.3, 63.6, 60.7, 58.9, 58.6, 59.9, 62.7, 66.9, 72, 77.6, 83.1)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 2).Value = $row[$c]
    }
}
